# Updates the "想去人数" (F column) counts on the "展览" and "全部类型" sheets,
# mirroring the refreshed numbers from the site data regeneration.

$wb = $excel.ActiveWorkbook

# New F-column values keyed by row number (same updates on both sheets).
$updates = @{
    2  = 1194
    3  = 604
    5  = 37
    6  = 193
    7  = 64
    8  = 65
    10 = 5565
    11 = 4943
    12 = 21
    13 = 45
    16 = 204
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
